$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Update the "Last updated:" date from "January 7, 2025" to
#    "April 29, 2025". The original paragraph is split across three
#    runs: "January 7" | ", 202" | "5" (all sharing identical rPr).
#    A plain text replace lets the engine merge same-formatted runs
#    together; toggling bold on/off right after editing the first
#    run forces Word to keep that run's boundary distinct from the
#    preceding "Last updated: " run, matching the target structure:
#       "Last updated: " | "April 29" | ", 2025"
# -----------------------------------------------------------------
$dateRange = $d.Content
$dateRange.Find.Execute("January 7", $true, $false, $false, $false, $false, $true, 1, $false, "April 29", 2) | Out-Null
$dateRange.Font.Bold = $true
$dateRange.Font.Bold = $false

# -----------------------------------------------------------------
# 2) Add the references/ChatGPT requirement. There are two empty
#    paragraphs right after "The conclusions section can include
#    future work, if there was more time." -- the first of those
#    gets the new sentence, and a brand-new empty paragraph (with
#    the same formatting) is inserted right after it.
# -----------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
  $p = $d.Paragraphs($i)
  if ($p.Range.Text -like "*The conclusions section can include future work*") {
    $target = $d.Paragraphs($i + 1)
    $target.Range.Text = "v. References. If you" + [char]0x2019 + "ve used tools like ChatGPT, mention how they were used."
    $target.Range.InsertParagraphAfter()
    break
  }
}

# -----------------------------------------------------------------
# 3) Merge "7. Team" + "work" (two runs, the second bold) into a
#    single bold run "7. Teamwork". A no-op find/replace of the
#    first run's text causes the adjoining identically-formatted
#    run to be absorbed into it.
# -----------------------------------------------------------------
$d.Content.Find.Execute("7. Team", $true, $false, $false, $false, $false, $true, 1, $false, "7. Team", 2) | Out-Null
